$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Flag" column (F3:F10) from "No" to "Yes"
$ws.Range("F3:F10").Value = "Yes"

# Update the active selection to F2:F10 (active cell F2)
$ws.Range("F2:F10").Select()
